# Generate Report for Handback
# Updates file identifiers / timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# New identifiers/filenames replacing the old ones.
$oldGuid1 = "185f8477-1d82-457c-9cf6-d4b222813430"
$newGuid1 = "5daf2757-4342-4441-8a1b-7c039d6e4905"
$oldGuid2 = "f06f1f58-d76a-490e-9da1-1e75b1249941"
$newGuid2 = "ffffbe4f7217-6524-4177-bc81-4f7a0057af32"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-08-25 03:02:33"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-08-25 03:02:33"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\$newGuid1.md"
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = "e2e\$newGuid2.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid1.md"
$wsZh.Range("I2").Value = "$newGuid1.md"
$wsZh.Range("G2").Value = "$newGuid1.d0ea2da892cb80347cf45bf2890c2371202062f7.zh-cn.xlf"
$wsZh.Range("J2").Value = "$newGuid1.d0ea2da892cb80347cf45bf2890c2371202062f7.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-25 03:02:28"
$wsZh.Range("K2").Value = "2016-08-25 03:02:46"

$wsZh.Range("A3").Value = "$newGuid2.md"
$wsZh.Range("I3").Value = "$newGuid2.md"
$wsZh.Range("G3").Value = "$newGuid1.d0ea2da892cb80347cf45bf2890c2371202062f7.de-de.xlf"
$wsZh.Range("J3").Value = "$newGuid1.d0ea2da892cb80347cf45bf2890c2371202062f7.de-de.xlf"
$wsZh.Range("H3").Value = "2016-08-25 03:02:28"
$wsZh.Range("K3").Value = "2016-08-25 03:02:46"

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "$newGuid2.md"
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = "$newGuid2.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid1.md"
$wsDe.Range("I2").Value = "$newGuid1.md"
$wsDe.Range("G2").Value = "$newGuid1.d0ea2da892cb80347cf45bf2890c2371202062f7.de-de.xlf"
$wsDe.Range("J2").Value = "$newGuid1.d0ea2da892cb80347cf45bf2890c2371202062f7.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-25 03:02:33"
$wsDe.Range("K2").Value = "2016-08-25 03:02:53"

$wsDe.Range("A3").Value = "$newGuid2.md"
$wsDe.Range("I3").Value = "$newGuid2.md"
$wsDe.Range("G3").Value = "$newGuid1.d0ea2da892cb80347cf45bf2890c2371202062f7.de-de.xlf"
$wsDe.Range("J3").Value = "$newGuid1.d0ea2da892cb80347cf45bf2890c2371202062f7.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-25 03:02:33"
$wsDe.Range("K3").Value = "2016-08-25 03:02:53"

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = "$newGuid1.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "$newGuid2.md"
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = "$newGuid2.md"
    }
}
